# Update the PUBCON YoY AVERAGE(1,9) forecast vector sheet: bugfixed evaluation
# and simulated rt_data for components. Every forecast row (date_of_forecast,
# y_0, y_0_forecast, y_1, y_1_forecast) shifts/changes and one new trailing row
# (row 53) is appended, extending the sheet from A1:E52 to A1:E53.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 53 is brand new. Copy the date-column formatting (style "s=2": custom
# YYYY-MM-DD HH:MM:SS number format + bold/centered/bordered style) from the
# previous last row (A52) down into A53 before writing its value.
$ws.Range("A52").Copy($ws.Range("A53"))

$ws.Cells.Item(2, 1).Value = 39400
$ws.Cells.Item(2, 2).Value = 2007
$ws.Cells.Item(2, 3).Value = 1.75539628881467
$ws.Cells.Item(2, 4).Value = 2008
$ws.Cells.Item(2, 5).Value = 0.2337905658324813

$ws.Cells.Item(3, 1).Value = 39583
$ws.Cells.Item(3, 2).Value = 2008
$ws.Cells.Item(3, 3).Value = 2.071001150990881
$ws.Cells.Item(3, 4).Value = 2009
$ws.Cells.Item(3, 5).Value = 5.125398993052044

$ws.Cells.Item(4, 1).Value = 39765
$ws.Cells.Item(4, 2).Value = 2008
$ws.Cells.Item(4, 3).Value = 2.213911448916162
$ws.Cells.Item(4, 4).Value = 2009
$ws.Cells.Item(4, 5).Value = 3.386383090739953

$ws.Cells.Item(5, 1).Value = 39948
$ws.Cells.Item(5, 2).Value = 2009
$ws.Cells.Item(5, 3).Value = 0.8787496612562951
$ws.Cells.Item(5, 4).Value = 2010
$ws.Cells.Item(5, 5).Value = 1.013823151053095

$ws.Cells.Item(6, 1).Value = 40130
$ws.Cells.Item(6, 2).Value = 2009
$ws.Cells.Item(6, 3).Value = 2.533533936850563
$ws.Cells.Item(6, 4).Value = 2010
$ws.Cells.Item(6, 5).Value = 0.984293482975751

$ws.Cells.Item(7, 1).Value = 40310
$ws.Cells.Item(7, 2).Value = 2010
$ws.Cells.Item(7, 3).Value = 2.723861837742825
$ws.Cells.Item(7, 4).Value = 2011
$ws.Cells.Item(7, 5).Value = 4.356912452939476

$ws.Cells.Item(8, 1).Value = 40494
$ws.Cells.Item(8, 2).Value = 2010
$ws.Cells.Item(8, 3).Value = 2.088987486264915
$ws.Cells.Item(8, 4).Value = 2011
$ws.Cells.Item(8, 5).Value = 3.612753212925401

$ws.Cells.Item(9, 1).Value = 40676
$ws.Cells.Item(9, 2).Value = 2011
$ws.Cells.Item(9, 3).Value = 2.431458940167008
$ws.Cells.Item(9, 4).Value = 2012
$ws.Cells.Item(9, 5).Value = 5.259925231829898

$ws.Cells.Item(10, 1).Value = 40862
$ws.Cells.Item(10, 2).Value = 2011
$ws.Cells.Item(10, 3).Value = 1.212544822741002
$ws.Cells.Item(10, 4).Value = 2012
$ws.Cells.Item(10, 5).Value = 2.158838189283174

$ws.Cells.Item(11, 1).Value = 41044
$ws.Cells.Item(11, 2).Value = 2012
$ws.Cells.Item(11, 3).Value = 1.447930496829541
$ws.Cells.Item(11, 4).Value = 2013
$ws.Cells.Item(11, 5).Value = 0.7814198158142105

$ws.Cells.Item(12, 1).Value = 41228
$ws.Cells.Item(12, 2).Value = 2012
$ws.Cells.Item(12, 3).Value = 1.196776590518644
$ws.Cells.Item(12, 4).Value = 2013
$ws.Cells.Item(12, 5).Value = 1.194058515117336

$ws.Cells.Item(13, 1).Value = 41409
$ws.Cells.Item(13, 2).Value = 2013
$ws.Cells.Item(13, 3).Value = 0.3494148569448852
$ws.Cells.Item(13, 4).Value = 2014
$ws.Cells.Item(13, 5).Value = -0.5632089058212553

$ws.Cells.Item(14, 1).Value = 41592
$ws.Cells.Item(14, 2).Value = 2013
$ws.Cells.Item(14, 3).Value = 0.4712609263772594
$ws.Cells.Item(14, 4).Value = 2014
$ws.Cells.Item(14, 5).Value = 1.409662779709819

$ws.Cells.Item(15, 1).Value = 41774
$ws.Cells.Item(15, 2).Value = 2014
$ws.Cells.Item(15, 3).Value = 0.6742451383205061
$ws.Cells.Item(15, 4).Value = 2015
$ws.Cells.Item(15, 5).Value = 1.713290556413605

$ws.Cells.Item(16, 1).Value = 41957
$ws.Cells.Item(16, 2).Value = 2014
$ws.Cells.Item(16, 3).Value = 0.8783377572271434
$ws.Cells.Item(16, 4).Value = 2015
$ws.Cells.Item(16, 5).Value = 2.372074663906587

$ws.Cells.Item(17, 1).Value = 42137
$ws.Cells.Item(17, 2).Value = 2015
$ws.Cells.Item(17, 3).Value = 2.178094576990031
$ws.Cells.Item(17, 4).Value = 2016
$ws.Cells.Item(17, 5).Value = 2.743122680804988

$ws.Cells.Item(18, 1).Value = 42321
$ws.Cells.Item(18, 2).Value = 2015
$ws.Cells.Item(18, 3).Value = 2.29066283401107
$ws.Cells.Item(18, 4).Value = 2016
$ws.Cells.Item(18, 5).Value = 4.595879021798321

$ws.Cells.Item(19, 1).Value = 42503
$ws.Cells.Item(19, 2).Value = 2016
$ws.Cells.Item(19, 3).Value = 2.623600596229347
$ws.Cells.Item(19, 4).Value = 2017
$ws.Cells.Item(19, 5).Value = 1.891220645890002

$ws.Cells.Item(20, 1).Value = 42689
$ws.Cells.Item(20, 2).Value = 2016
$ws.Cells.Item(20, 3).Value = 4.109890522944348
$ws.Cells.Item(20, 4).Value = 2017
$ws.Cells.Item(20, 5).Value = 4.034919509273061

$ws.Cells.Item(21, 1).Value = 42867
$ws.Cells.Item(21, 2).Value = 2017
$ws.Cells.Item(21, 3).Value = 1.715791310593229
$ws.Cells.Item(21, 4).Value = 2018
$ws.Cells.Item(21, 5).Value = 1.687339605296523

$ws.Cells.Item(22, 1).Value = 43053
$ws.Cells.Item(22, 2).Value = 2017
$ws.Cells.Item(22, 3).Value = 1.336316831462692
$ws.Cells.Item(22, 4).Value = 2018
$ws.Cells.Item(22, 5).Value = 0.02883756256675252

$ws.Cells.Item(23, 1).Value = 43145
$ws.Cells.Item(23, 2).Value = 2018
$ws.Cells.Item(23, 3).Value = 1.862609889357314
$ws.Cells.Item(23, 4).Value = 2019
$ws.Cells.Item(23, 5).Value = 2.016108144125295

$ws.Cells.Item(24, 1).Value = 43235
$ws.Cells.Item(24, 2).Value = 2018
$ws.Cells.Item(24, 3).Value = 0.3317798769387315
$ws.Cells.Item(24, 4).Value = 2019
$ws.Cells.Item(24, 5).Value = -2.079848588862154

$ws.Cells.Item(25, 1).Value = 43326
$ws.Cells.Item(25, 2).Value = 2018
$ws.Cells.Item(25, 3).Value = 1.369357173039498
$ws.Cells.Item(25, 4).Value = 2019
$ws.Cells.Item(25, 5).Value = 2.147322685428343

$ws.Cells.Item(26, 1).Value = 43418
$ws.Cells.Item(26, 2).Value = 2018
$ws.Cells.Item(26, 3).Value = 1.197912858979611
$ws.Cells.Item(26, 4).Value = 2019
$ws.Cells.Item(26, 5).Value = 0.9262553939922924

$ws.Cells.Item(27, 1).Value = 43510
$ws.Cells.Item(27, 2).Value = 2019
$ws.Cells.Item(27, 3).Value = 4.206171608278875
$ws.Cells.Item(27, 4).Value = 2020
$ws.Cells.Item(27, 5).Value = 6.673057635354218

$ws.Cells.Item(28, 1).Value = 43600
$ws.Cells.Item(28, 2).Value = 2019
$ws.Cells.Item(28, 3).Value = 0.8934982674867697
$ws.Cells.Item(28, 4).Value = 2020
$ws.Cells.Item(28, 5).Value = -1.194610791899997

$ws.Cells.Item(29, 1).Value = 43691
$ws.Cells.Item(29, 2).Value = 2019
$ws.Cells.Item(29, 3).Value = 1.749048192229496
$ws.Cells.Item(29, 4).Value = 2020
$ws.Cells.Item(29, 5).Value = 2.097071260293459

$ws.Cells.Item(30, 1).Value = 43783
$ws.Cells.Item(30, 2).Value = 2019
$ws.Cells.Item(30, 3).Value = 1.727537197898665
$ws.Cells.Item(30, 4).Value = 2020
$ws.Cells.Item(30, 5).Value = 2.928189816005666

$ws.Cells.Item(31, 1).Value = 43875
$ws.Cells.Item(31, 2).Value = 2020
$ws.Cells.Item(31, 3).Value = 2.306826470345347
$ws.Cells.Item(31, 4).Value = 2021
$ws.Cells.Item(31, 5).Value = 1.405199465409468

$ws.Cells.Item(32, 1).Value = 43966
$ws.Cells.Item(32, 2).Value = 2020
$ws.Cells.Item(32, 3).Value = 1.979074033580819
$ws.Cells.Item(32, 4).Value = 2021
$ws.Cells.Item(32, 5).Value = 0.8024032015999882

$ws.Cells.Item(33, 1).Value = 44068
$ws.Cells.Item(33, 2).Value = 2020
$ws.Cells.Item(33, 3).Value = 2.980209378995857
$ws.Cells.Item(33, 4).Value = 2021
$ws.Cells.Item(33, 5).Value = 2.026374749120596

$ws.Cells.Item(34, 1).Value = 44159
$ws.Cells.Item(34, 2).Value = 2020
$ws.Cells.Item(34, 3).Value = 3.647228437274408
$ws.Cells.Item(34, 4).Value = 2021
$ws.Cells.Item(34, 5).Value = 3.673004547855219

$ws.Cells.Item(35, 1).Value = 44251
$ws.Cells.Item(35, 2).Value = 2021
$ws.Cells.Item(35, 3).Value = 2.128447415063373
$ws.Cells.Item(35, 4).Value = 2022
$ws.Cells.Item(35, 5).Value = 2.829537440100038

$ws.Cells.Item(36, 1).Value = 44341
$ws.Cells.Item(36, 2).Value = 2021
$ws.Cells.Item(36, 3).Value = 2.552476296061434
$ws.Cells.Item(36, 4).Value = 2022
$ws.Cells.Item(36, 5).Value = 3.086122033237126

$ws.Cells.Item(37, 1).Value = 44432
$ws.Cells.Item(37, 2).Value = 2021
$ws.Cells.Item(37, 3).Value = 1.954146674711188
$ws.Cells.Item(37, 4).Value = 2022
$ws.Cells.Item(37, 5).Value = -1.871731962523027

$ws.Cells.Item(38, 1).Value = 44525
$ws.Cells.Item(38, 2).Value = 2021
$ws.Cells.Item(38, 3).Value = 2.777797690741424
$ws.Cells.Item(38, 4).Value = 2022
$ws.Cells.Item(38, 5).Value = 1.579011422502852

$ws.Cells.Item(39, 1).Value = 44617
$ws.Cells.Item(39, 2).Value = 2022
$ws.Cells.Item(39, 3).Value = 1.014265466411501
$ws.Cells.Item(39, 4).Value = 2023
$ws.Cells.Item(39, 5).Value = 0.00000000000006661338147750939

$ws.Cells.Item(40, 1).Value = 44706
$ws.Cells.Item(40, 2).Value = 2022
$ws.Cells.Item(40, 3).Value = 0.3901728183783204
$ws.Cells.Item(40, 4).Value = 2023
$ws.Cells.Item(40, 5).Value = 1.906002353653125

$ws.Cells.Item(41, 1).Value = 44798
$ws.Cells.Item(41, 2).Value = 2022
$ws.Cells.Item(41, 3).Value = 2.69102598245059
$ws.Cells.Item(41, 4).Value = 2023
$ws.Cells.Item(41, 5).Value = 3.932998599265303

$ws.Cells.Item(42, 1).Value = 44890
$ws.Cells.Item(42, 2).Value = 2022
$ws.Cells.Item(42, 3).Value = 0.6994919452575576
$ws.Cells.Item(42, 4).Value = 2023
$ws.Cells.Item(42, 5).Value = -2.087978868409623

$ws.Cells.Item(43, 1).Value = 44981
$ws.Cells.Item(43, 2).Value = 2023
$ws.Cells.Item(43, 3).Value = -0.6123001687638907
$ws.Cells.Item(43, 4).Value = 2024
$ws.Cells.Item(43, 5).Value = -1.985049937500016

$ws.Cells.Item(44, 1).Value = 45071
$ws.Cells.Item(44, 2).Value = 2023
$ws.Cells.Item(44, 3).Value = -1.906744368254853
$ws.Cells.Item(44, 4).Value = 2024
$ws.Cells.Item(44, 5).Value = 9.556389850000446

$ws.Cells.Item(45, 1).Value = 45163
$ws.Cells.Item(45, 2).Value = 2023
$ws.Cells.Item(45, 3).Value = -1.669605379075589
$ws.Cells.Item(45, 4).Value = 2024
$ws.Cells.Item(45, 5).Value = -0.4516362914963601

$ws.Cells.Item(46, 1).Value = 45254
$ws.Cells.Item(46, 2).Value = 2023
$ws.Cells.Item(46, 3).Value = -1.432689847121871
$ws.Cells.Item(46, 4).Value = 2024
$ws.Cells.Item(46, 5).Value = 0.1172571542027212

$ws.Cells.Item(47, 1).Value = 45345
$ws.Cells.Item(47, 2).Value = 2024
$ws.Cells.Item(47, 3).Value = 0.5893837960974757
$ws.Cells.Item(47, 4).Value = 2025
$ws.Cells.Item(47, 5).Value = -0.3994003998999518

$ws.Cells.Item(48, 1).Value = 45436
$ws.Cells.Item(48, 2).Value = 2024
$ws.Cells.Item(48, 3).Value = 1.152760694685062
$ws.Cells.Item(48, 4).Value = 2025
$ws.Cells.Item(48, 5).Value = 0.6444718444275521

$ws.Cells.Item(49, 1).Value = 45534
$ws.Cells.Item(49, 2).Value = 2024
$ws.Cells.Item(49, 3).Value = 1.780300968358017
$ws.Cells.Item(49, 4).Value = 2025
$ws.Cells.Item(49, 5).Value = 0.8104204674762139

$ws.Cells.Item(50, 1).Value = 45618
$ws.Cells.Item(50, 2).Value = 2024
$ws.Cells.Item(50, 3).Value = 2.033479419175133
$ws.Cells.Item(50, 4).Value = 2025
$ws.Cells.Item(50, 5).Value = 1.317145539573517

$ws.Cells.Item(51, 1).Value = 45713
$ws.Cells.Item(51, 2).Value = 2025
$ws.Cells.Item(51, 3).Value = 2.506151357112452
$ws.Cells.Item(51, 4).Value = 2026
$ws.Cells.Item(51, 5).Value = 0

$ws.Cells.Item(52, 1).Value = 45800
$ws.Cells.Item(52, 2).Value = 2025
$ws.Cells.Item(52, 3).Value = 2.287212358310953
$ws.Cells.Item(52, 4).Value = 2026
$ws.Cells.Item(52, 5).Value = 1.801540135156521

$ws.Cells.Item(53, 1).Value = 45891
$ws.Cells.Item(53, 2).Value = 2025
$ws.Cells.Item(53, 3).Value = 2.481068287768839
$ws.Cells.Item(53, 4).Value = 2026
$ws.Cells.Item(53, 5).Value = 1.642273543406181
